# Update the cryptos list snapshot on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Coin, Link, Price, Volume(1h))
$rows = @{
    2  = @('Bitcoin',                        'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc',                 '67.484.83',   '  -0.94%  ')
    3  = @('Ethereum',                       'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth',                '3.312.27',    '  +1.25%  ')
    4  = @('TetherUSD',                      'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt',              '1.00',        '  -0.03%  ')
    5  = @('Solana',                         'https://coinranking.com/coin/zNZHO_Sjf+solana-sol',                      '186.66',      '  +1.17%  ')
    6  = @('BNB',                            'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb',                     '578.37',      '  -0.92%  ')
    7  = @('USDC',                           'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc',                   '1.00',        '  +0.06%  ')
    8  = @('XRP',                            'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp',                     '0.606',       '  +0.45%  ')
    9  = @('Dogecoin',                       'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge',               '0.130',       '  -0.24%  ')
    10 = @('Toncoin',                        'https://coinranking.com/coin/67YlI0K1b+toncoin-ton',                     '6.68',        '  +1.08%  ')
    11 = @('Cardano',                        'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada',                 '0.409',       '  -0.27%  ')
    12 = @('WrappedliquidstakedEther2.0',     'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','3.887.91',    '  +1.27%  ')
    13 = @('TRON',                           'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx',                    '0.138',       '  -0.38%  ')
    14 = @('Avalanche',                      'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax',                  '27.51',       '  -0.06%  ')
    15 = @('WrappedBTC',                     'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc',              '67.723.48',   '  -0.61%  ')
    16 = @('ShibaInu',                       'https://coinranking.com/coin/xz24e0BjL+shibainu-shib',                   '0.0000168',   '  -0.26%  ')
    17 = @('WrappedEther',                   'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',          '3.295.62',    '  +0.93%  ')
    18 = @('BitcoinCash',                    'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',             '445.76',      '  +6.65%  ')
    19 = @('Polkadot',                       'https://coinranking.com/coin/25W7FG7om+polkadot-dot',                    '5.71',        '  -0.81%  ')
    20 = @('Chainlink',                      'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link',              '13.62',       '  +1.95%  ')
    21 = @('Uniswap',                        'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni',                     '7.76',        '  +2.65%  ')
    22 = @('Litecoin',                       'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc',                '74.06',       '  +3.71%  ')
    23 = @('Dai',                            'https://coinranking.com/coin/MoTuySvg7+dai-dai',                         '1.00',        '  -0.12%  ')
    24 = @('Polygon',                        'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic',                '0.519',       '  +2.10%  ')
    25 = @('WrappedeETH',                    'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth',               '3.455.34',    '  +1.22%  ')
    26 = @('PEPE',                           'https://coinranking.com/coin/03WI8NQPF+pepe-pepe',                       '0.0000119',   '  +1.38%  ')
    27 = @('Kaspa',                          'https://coinranking.com/coin/V8GxkwWow+kaspa-kas',                       '0.189',       '  +1.08%  ')
    28 = @('InternetComputer(DFINITY)',      'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp',     '9.09',        '  -3.60%  ')
    29 = @('Binance-PegBSC-USD',             'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd',      '1.00',        '  +0.14%  ')
    30 = @('PancakeSwap',                    'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake',                '1.98',        '  +1.53%  ')
    31 = @('EthereumClassic',                'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc',         '22.99',       '  +1.14%  ')
    32 = @('NEARProtocol',                   'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',               '5.35',        '  -2.24%  ')
    33 = @('Fetch.AI',                       'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet',                 '1.25',        '  +0.14%  ')
    34 = @('USDe',                           'https://coinranking.com/coin/exbfr2U-0+usde-usde',                       '0.999',       '  -0.01%  ')
    35 = @('Aptos',                          'https://coinranking.com/coin/HGYj5JCv5+aptos-apt',                       '6.82',        '  -1.16%  ')
    36 = @('ImmutableX',                     'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',                  '1.52',        '  +4.98%  ')
    37 = @('Monero',                         'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',                  '163.06',      '  -0.49%  ')
    38 = @('Stacks',                         'https://coinranking.com/coin/mMPrMcB7+stacks-stx',                       '1.86',        '  -1.61%  ')
    39 = @('EnergySwap',                     'https://coinranking.com/coin/SbWqqTui-+energyswap-ens',                  '27.23',       '  +0.12%  ')
    40 = @('Mantle',                         'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt',                      '0.792',       '  -0.89%  ')
    41 = @('Filecoin',                       'https://coinranking.com/coin/ymQub4fuB+filecoin-fil',                    '4.48',        '  +0.22%  ')
    42 = @('Maker',                          'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr',                   '2.760.43',    '  +3.57%  ')
    43 = @('RenderToken',                    'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr',            '6.27',        '  -0.81%  ')
    44 = @('dogwifhat',                      'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif',                   '2.43',        '  -0.30%  ')
    45 = @('InjectiveProtocol',              'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj',           '24.97',       '  +1.88%  ')
    46 = @('Hedera',                         'https://coinranking.com/coin/jad286TjB+hedera-hbar',                     '0.0676',      '  -0.36%  ')
    47 = @('OKB',                            'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb',                     '40.24',       '  -1.65%  ')
    48 = @('Bittensor',                      'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao',                  '327.20',      '  -3.29%  ')
    49 = @('VeChain',                        'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',                 '0.0275',      '  +0.27%  ')
    50 = @('ONDO',                           'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo',                       '0.995',       '  +1.62%  ')
    51 = @('Arweave',                        'https://coinranking.com/coin/7XWg41D1+arweave-ar',                       '31.31',       '  +1.82%  ')
}

# Capture the plain/default cell style (used by all data cells B2:E51 originally,
# i.e. no special formatting) so we can restore it after forcing text entry.
$defaultStyle = $ws.Cells.Item(2, 2).Style

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $cellE = $ws.Cells.Item($r, 5)

    # Temporarily force a text number format so numeric-looking strings
    # (e.g. "67.484.83", "1.00", "0.0000168") are preserved verbatim as text
    # instead of being parsed/rounded as numbers, then restore the original
    # (default/general) style so the saved file's formatting is unchanged.
    $cellB.NumberFormat = "@"
    $cellB.Value = $vals[0]
    $cellB.Style = $defaultStyle

    $cellC.NumberFormat = "@"
    $cellC.Value = $vals[1]
    $cellC.Style = $defaultStyle

    $cellD.NumberFormat = "@"
    $cellD.Value = $vals[2]
    $cellD.Style = $defaultStyle

    $cellE.NumberFormat = "@"
    $cellE.Value = $vals[3]
    $cellE.Style = $defaultStyle
}
